# Final code for deriving SoA and RTs from the data!
#
# Inserts 9 new "reaction time" (RT) columns right after column G (neut_r)
# and before the old column H (comp), shifting every existing column from
# H:AU out to Q:BD. The new columns get header labels in row 1 and the new
# per-subject RT values in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 blank columns at H:P, shifting existing H:AU content (and its
# formatting) right to Q:BD. This preserves every pre-existing value /
# shared-string reference exactly - only its column position changes.
$ws.Range("H1:P2").Insert(-4161) | Out-Null   # xlShiftToRight

# The inserted header cells (H1:P1) already carry the bold/centered/bordered
# header style (it rides along with the insert), so we only need to fill in
# the text.
$ws.Range("H1").Value = "comp_r_RT"
$ws.Range("I1").Value = "comp_l_RT"
$ws.Range("J1").Value = "incomp_r_RT"
$ws.Range("K1").Value = "incomp_l_RT"
$ws.Range("L1").Value = "neut_r_RT"
$ws.Range("M1").Value = "neut_l_RT"
$ws.Range("N1").Value = "comp_RT"
$ws.Range("O1").Value = "incomp_RT"
$ws.Range("P1").Value = "neut_RT"

# New reaction-time data for row 2.
$ws.Range("H2").Value = 319.426666666667
$ws.Range("I2").Value = 326.554054054054
$ws.Range("J2").Value = 352.848484848485
$ws.Range("K2").Value = 341.161764705882
$ws.Range("L2").Value = 351.088235294118
$ws.Range("M2").Value = 337.376811594203
$ws.Range("N2").Value = 322.99036036036
$ws.Range("O2").Value = 347.005124777184
$ws.Range("P2").Value = 344.23252344416

"done"
